# MOSIP-17570 added supervisor rejected email and sms templates
#
# Adds 18 new rows (1726-1743) to Sheet1 describing three new template
# types - RPR_SUP_REJECT_EMAIL, RPR_SUP_REJECT_SMS and
# RPR_SUP_REJECT_EMAIL_SUBJECT - replicated across the six languages
# already present in the sheet (eng, fra, ara, hin, kan, tam).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# lang_code, code, descr triples used to build the new rows. The
# is_active (column D) value is copied from an existing row further
# down so the new cells keep the same shared-string/style encoding
# ("TRUE" as text, not boolean) as the rest of the sheet.
$newRows = @(
    @("eng", "RPR_SUP_REJECT_EMAIL", "Template for Supervisor Reject Email"),
    @("eng", "RPR_SUP_REJECT_SMS", "Template for Supervisor Reject SMS"),
    @("fra", "RPR_SUP_REJECT_EMAIL", "Template for Supervisor Reject Email"),
    @("fra", "RPR_SUP_REJECT_SMS", "Template for Supervisor Reject SMS"),
    @("ara", "RPR_SUP_REJECT_EMAIL", "Template for Supervisor Reject Email"),
    @("ara", "RPR_SUP_REJECT_SMS", "Template for Supervisor Reject SMS"),
    @("hin", "RPR_SUP_REJECT_EMAIL", "Template for Supervisor Reject Email"),
    @("hin", "RPR_SUP_REJECT_SMS", "Template for Supervisor Reject SMS"),
    @("kan", "RPR_SUP_REJECT_EMAIL", "Template for Supervisor Reject Email"),
    @("kan", "RPR_SUP_REJECT_SMS", "Template for Supervisor Reject SMS"),
    @("tam", "RPR_SUP_REJECT_EMAIL", "Template for Supervisor Reject Email"),
    @("tam", "RPR_SUP_REJECT_SMS", "Template for Supervisor Reject SMS"),
    @("eng", "RPR_SUP_REJECT_EMAIL_SUBJECT", "Template for Supervisor Reject Email Subject"),
    @("fra", "RPR_SUP_REJECT_EMAIL_SUBJECT", "Template for Supervisor Reject Email Subject"),
    @("ara", "RPR_SUP_REJECT_EMAIL_SUBJECT", "Template for Supervisor Reject Email Subject"),
    @("hin", "RPR_SUP_REJECT_EMAIL_SUBJECT", "Template for Supervisor Reject Email Subject"),
    @("kan", "RPR_SUP_REJECT_EMAIL_SUBJECT", "Template for Supervisor Reject Email Subject"),
    @("tam", "RPR_SUP_REJECT_EMAIL_SUBJECT", "Template for Supervisor Reject Email Subject")
)

$startRow = 1726
$row = $startRow
foreach ($item in $newRows) {
    $ws.Range("A$row").Value = $item[0]
    $ws.Range("B$row").Value = $item[1]
    $ws.Range("C$row").Value = $item[2]

    # Column D ("is_active") is always the text "TRUE" elsewhere in the
    # sheet; copy an existing cell so the new one matches exactly
    # instead of Excel auto-coercing a literal "TRUE" into a boolean.
    $ws.Range("D2").Copy()
    $ws.Range("D$row").PasteSpecial()

    $row = $row + 1
}

$excel.CutCopyMode = $false

# Update the view to match where the author was working when the rows
# were added.
$ws.Range("F1730").Select() | Out-Null
